# --- Sheet1: add "Invoice Numer" as a new first column ---
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Columns.Item(1).Insert()

$ws1.Range("A1").Value = "Invoice Numer"
$ws1.Range("B1").Copy()
$ws1.Range("A1").PasteSpecial(-4122)
$ws1.Columns.Item(1).ColumnWidth = 13.25

$ws1.Range("A2").Value = 1000001
$ws1.Range("A3").Value = 1000002
$ws1.Range("A4").Value = 1000003
$ws1.Range("A5").Value = 1000004
$ws1.Range("A6").Value = 1000005
$ws1.Range("A7").Value = 1000006
$ws1.Range("A8").Value = 1000007
$ws1.Range("A9").Value = 1000008
$ws1.Range("A10").Value = 1000009
$ws1.Range("A11").Value = 1000010
$ws1.Range("A12").Value = 1000011
$ws1.Range("A13").Value = 1000012
$ws1.Range("A14").Value = 1000013

$ws1.Range("B21").Select()

# --- Summary sheet: rename "Summary" -> "T5 Summary", add "BART Summary" column ---
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B1").Value = "T5 Summary"
$ws2.Range("B2").Value = "there is cash crunch impacting fund transfer. partial payment expected by Jan Month end with further settlements in Feb. Client has released full payment for all open invoices. expect all payments to be applied by 5th April in the system."
$ws2.Range("B3").Value = "B LLC has expressed an inability to pay at the moment and promised to make a payment by 1st march. unable to commit to a new date, but instead have mentioned payment in ""near future"" the client will be sent to bad debt collection; payment is not expected and might need a write-off."

$ws2.Range("C1").Value = "BART Summary"
$ws2.Range("A1").Copy()
$ws2.Range("C1").PasteSpecial(-4122)

$ws2.Range("C2").Value = "Update On: 15th Mar 2024, AP Contact has been changed from John to Matthew effective immediately. Matthew will be the SPOC for all payments going forward. Expect all payments to be applied by 5th April in the system."
$ws2.Range("C3").Value = "Update On: 13th Jan 2024, Sent initial chaser to client on outstanding balance. Per email from Jim (AP), B LLC has expressed an inability to pay at the moment and promised to make a payment by 1st March."
